{"js": "const replacements = [\n  [\"183\u00f79=\", \"843\u00f73=\"],\n  [\"147\u00f74=\", \"935\u00f75=\"],\n  [\"530\u00f79=\", \"544\u00f79=\"],\n  [\"591\u00f74=\", \"411\u00f72=\"],\n  [\"287\u00f75=\", \"583\u00f77=\"],\n  [\"719\u00f72=\", \"874\u00f74=\"],\n  [\"922\u00f75=\", \"210\u00f77=\"],\n  [\"274\u00f72=\", \"505\u00f76=\"],\n  [\"645\u00f72=\", \"261\u00f77=\"],\n  [\"226\u00f78=\", \"536\u00f78=\"],\n  [\"284\u00f77=\", \"534\u00f76=\"],\n  [\"911\u00f74=\", \"601\u00f76=\"],\n  [\"744\u00f75=\", \"166\u00f77=\"],\n  [\"903\u00f74=\", \"699\u00f72=\"],\n  [\"342\u00f79=\", \"525\u00f72=\"],\n  [\"114\u00f76=\", \"272\u00f72=\"],\n  [\"666\u00f74=\", \"900\u00f73=\"],\n  [\"367\u00f78=\", \"255\u00f73=\"],\n  [\"970\u00f77=\", \"511\u00f76=\"],\n  [\"854\u00f79=\", \"504\u00f73=\"],\n  [\"920\u00f73=\", \"142\u00f78=\"],\n  [\"995\u00f75=\", \"501\u00f78=\"],\n  [\"216\u00f73=\", \"629\u00f75=\"],\n  [\"245\u00f74=\", \"229\u00f72=\"],\n  [\"329\u00f75=\", \"647\u00f75=\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"183\u00f79=\"; New=\"843\u00f73=\"},\n    @{Old=\"147\u00f74=\"; New=\"935\u00f75=\"},\n    @{Old=\"530\u00f79=\"; New=\"544\u00f79=\"},\n    @{Old=\"591\u00f74=\"; New=\"411\u00f72=\"},\n    @{Old=\"287\u00f75=\"; New=\"583\u00f77=\"},\n    @{Old=\"719\u00f72=\"; New=\"874\u00f74=\"},\n    @{Old=\"922\u00f75=\"; New=\"210\u00f77=\"},\n    @{Old=\"274\u00f72=\"; New=\"505\u00f76=\"},\n    @{Old=\"645\u00f72=\"; New=\"261\u00f77=\"},\n    @{Old=\"226\u00f78=\"; New=\"536\u00f78=\"},\n    @{Old=\"284\u00f77=\"; New=\"534\u00f76=\"},\n    @{Old=\"911\u00f74=\"; New=\"601\u00f76=\"},\n    @{Old=\"744\u00f75=\"; New=\"166\u00f77=\"},\n    @{Old=\"903\u00f74=\"; New=\"699\u00f72=\"},\n    @{Old=\"342\u00f79=\"; New=\"525\u00f72=\"},\n    @{Old=\"114\u00f76=\"; New=\"272\u00f72=\"},\n    @{Old=\"666\u00f74=\"; New=\"900\u00f73=\"},\n    @{Old=\"367\u00f78=\"; New=\"255\u00f73=\"},\n    @{Old=\"970\u00f77=\"; New=\"511\u00f76=\"},\n    @{Old=\"854\u00f79=\"; New=\"504\u00f73=\"},\n    @{Old=\"920\u00f73=\"; New=\"142\u00f78=\"},\n    @{Old=\"995\u00f75=\"; New=\"501\u00f78=\"},\n    @{Old=\"216\u00f73=\"; New=\"629\u00f75=\"},\n    @{Old=\"245\u00f74=\"; New=\"229\u00f72=\"},\n    @{Old=\"329\u00f75=\"; New=\"647\u00f75=\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair.Old, $false, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)\n}\n"}
